$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Total Annual Cost")

$updates = @(
    @{ Row = 2; Value = 42136288.29979742 },
    @{ Row = 3; Value = 41855920.77810784 },
    @{ Row = 5; Value = 41768311.54180558 },
    @{ Row = 6; Value = 41660848.67032823 },
    @{ Row = 7; Value = 41825962.86315875 },
    @{ Row = 8; Value = 41946959.37159097 },
    @{ Row = 10; Value = 41737656.85518473 },
    @{ Row = 16; Value = 42049247.41568523 },
    @{ Row = 17; Value = 41878548.21458574 },
    @{ Row = 18; Value = 42164494.02776385 },
    @{ Row = 20; Value = 41903771.35765467 },
    @{ Row = 21; Value = 41832669.08389357 },
    @{ Row = 22; Value = 41825558.85651746 },
    @{ Row = 23; Value = 41824847.83377985 },
    @{ Row = 25; Value = 41928042.1546158 },
    @{ Row = 26; Value = 41835096.16358968 },
    @{ Row = 28; Value = 41699113.73920813 },
    @{ Row = 29; Value = 41690838.25945077 },
    @{ Row = 30; Value = 41653101.34209275 },
    @{ Row = 31; Value = 41825188.1303352 },
    @{ Row = 32; Value = 41824810.76116162 },
    @{ Row = 33; Value = 41946844.16139126 },
    @{ Row = 34; Value = 41648908.77057428 },
    @{ Row = 35; Value = 41648908.39320511 },
    @{ Row = 36; Value = 41824768.83544643 },
    @{ Row = 38; Value = 41917121.64000935 },
    @{ Row = 39; Value = 41956075.24927603 },
    @{ Row = 40; Value = 41649831.87936275 },
    @{ Row = 41; Value = 41649000.70408395 },
    @{ Row = 42; Value = 41648917.58655607 },
    @{ Row = 43; Value = 41824769.75478153 },
    @{ Row = 45; Value = 41940883.91634111 },
    @{ Row = 46; Value = 41958451.47690921 },
    @{ Row = 47; Value = 41825929.98210432 },
    @{ Row = 48; Value = 41946956.08348553 },
    @{ Row = 49; Value = 41824780.44276195 },
    @{ Row = 50; Value = 41824769.9924043 },
    @{ Row = 51; Value = 41946840.08451553 },
    @{ Row = 52; Value = 41648908.3628867 },
    @{ Row = 53; Value = 41824768.8324146 },
    @{ Row = 54; Value = 41824768.83136956 },
    @{ Row = 55; Value = 41824768.83126505 },
    @{ Row = 58; Value = 42057344.82294449 },
    @{ Row = 59; Value = 41848026.43042255 },
    @{ Row = 60; Value = 41827094.59117036 },
    @{ Row = 61; Value = 41947072.54439213 },
    @{ Row = 62; Value = 41824792.08885261 },
    @{ Row = 63; Value = 41824771.15701336 },
    @{ Row = 64; Value = 41824769.06382944 },
    @{ Row = 65; Value = 41946839.99165804 },
    @{ Row = 66; Value = 41824768.8335792 },
    @{ Row = 68; Value = 41680631.57749583 },
    @{ Row = 69; Value = 41685997.94368104 },
    @{ Row = 73; Value = 42109855.58194485 },
    @{ Row = 75; Value = 41887058.60381541 },
    @{ Row = 76; Value = 41830997.80850964 },
    @{ Row = 77; Value = 41825391.72897907 },
    @{ Row = 78; Value = 41824831.12102601 },
    @{ Row = 79; Value = 41824775.0602307 },
    @{ Row = 80; Value = 41824769.45415117 },
    @{ Row = 82; Value = 41676964.40248317 },
    @{ Row = 83; Value = 41651713.95639599 },
    @{ Row = 84; Value = 41649188.91178727 },
    @{ Row = 85; Value = 41824796.88730466 },
    @{ Row = 86; Value = 41824771.63685857 },
    @{ Row = 87; Value = 41946840.24896096 },
    @{ Row = 88; Value = 41648908.37933125 },
    @{ Row = 90; Value = 41857471.83560449 },
    @{ Row = 91; Value = 41828039.13168855 },
    @{ Row = 92; Value = 41947166.99844395 },
    @{ Row = 93; Value = 41648941.05427954 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Value
}
